$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: replace the obsolete 4.7mF capacitor supplier part with the new one
$ws.Range("F23").Value = "667-EEU-FS0J472B"
$ws.Range("C23").Value = "Aluminum Electrolytic Capacitors - Radial Leaded 6.3VDC 4700uF 10000H 12.5x20mm "

# Row 35: previously empty, now populated with a new "1M Widerstand" line item
$ws.Range("B35").Value = "1M Widerstand"
$ws.Range("C35").Value = "Thick Film Resistors - SMD 1 MOhms 125 mW 0805 1%"
$ws.Range("E35").Value = "Mouser"
$ws.Range("F35").Value = "603-RC0805FR-071ML"
$ws.Rows.Item(35).RowHeight = 33

# Row 36: previously empty, now populated with a new "4.7nF Kondensator" line item
$ws.Range("F36").Value = "603-CC805KRX7R9BB472"
$ws.Range("C36").Value = "Multilayer Ceramic Capacitors MLCC - SMD/SMT 50V 4700pF X7R 0805 10%"
$ws.Range("E36").Value = "Mouser"
$ws.Range("B36").Value = "4.7nF Kondensator"
$ws.Rows.Item(36).RowHeight = 66

# Update the view state (scroll position / active selection) to reflect where the user left off
$ws.Application.ActiveWindow.ScrollRow = 26
$ws.Range("B37").Select()
